$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Sending/Ligand/Receptor/Target cluster labels and recompute metrics for rows 2-10
# (adds the previously-missing "ECs" sending-cluster rows per Dr Hou's advice)

# Row 2: ECs -> Lamc2 -> Itgb1 -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lamc2"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.4421816666666667
$ws.Cells.Item(2, 8).Value = 1.326545
$ws.Cells.Item(2, 9).Value = 0.06026482003168283
$ws.Cells.Item(2, 10).Value = 0.06026482003168283
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 117.044563
$ws.Cells.Item(2, 14).Value = 351.133689
$ws.Cells.Item(2, 15).Value = 0.3245365645427815
$ws.Cells.Item(2, 16).Value = 0.3245365645427815
$ws.Cells.Item(2, 17).Value = 51.75495994161167
$ws.Cells.Item(2, 18).Value = 465.7946394745051
$ws.Cells.Item(2, 19).Value = 0.01955813765587134
$ws.Cells.Item(2, 20).Value = 0.01955813765587134

# Row 3: ECs -> Lamc2 -> Itgb1 -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lamc2"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.4421816666666667
$ws.Cells.Item(3, 8).Value = 1.326545
$ws.Cells.Item(3, 9).Value = 0.06026482003168283
$ws.Cells.Item(3, 10).Value = 0.06026482003168283
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 101.5800373333333
$ws.Cells.Item(3, 14).Value = 304.740112
$ws.Cells.Item(3, 15).Value = 0.281657135515876
$ws.Cells.Item(3, 16).Value = 0.281657135515876
$ws.Cells.Item(3, 17).Value = 44.91683020811556
$ws.Cells.Item(3, 18).Value = 404.25147187304
$ws.Cells.Item(3, 19).Value = 0.01697401658250357
$ws.Cells.Item(3, 20).Value = 0.01697401658250356

# Row 4: ECs -> Lamc2 -> Itgb1 -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lamc2"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.4421816666666667
$ws.Cells.Item(4, 8).Value = 1.326545
$ws.Cells.Item(4, 9).Value = 0.06026482003168283
$ws.Cells.Item(4, 10).Value = 0.06026482003168283
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 142.0267893333333
$ws.Cells.Item(4, 14).Value = 426.080368
$ws.Cells.Item(4, 15).Value = 0.3938062999413425
$ws.Cells.Item(4, 16).Value = 0.3938062999413425
$ws.Cells.Item(4, 17).Value = 62.8016424187289
$ws.Cells.Item(4, 18).Value = 565.2147817685601
$ws.Cells.Item(4, 19).Value = 0.02373266579330792
$ws.Cells.Item(4, 20).Value = 0.02373266579330792

# Row 5: FAPs -> Lamc2 -> Itgb1 -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lamc2"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 6.369908666666666
$ws.Cells.Item(5, 8).Value = 19.109726
$ws.Cells.Item(5, 9).Value = 0.8681531333236113
$ws.Cells.Item(5, 10).Value = 0.8681531333236113
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 117.044563
$ws.Cells.Item(5, 14).Value = 351.133689
$ws.Cells.Item(5, 15).Value = 0.3245365645427815
$ws.Cells.Item(5, 16).Value = 0.3245365645427815
$ws.Cells.Item(5, 17).Value = 745.5631762399127
$ws.Cells.Item(5, 18).Value = 6710.068586159214
$ws.Cells.Item(5, 19).Value = 0.2817474353858961
$ws.Cells.Item(5, 20).Value = 0.2817474353858961

# Row 6: FAPs -> Lamc2 -> Itgb1 -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lamc2"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6.369908666666666
$ws.Cells.Item(6, 8).Value = 19.109726
$ws.Cells.Item(6, 9).Value = 0.8681531333236113
$ws.Cells.Item(6, 10).Value = 0.8681531333236113
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 101.5800373333333
$ws.Cells.Item(6, 14).Value = 304.740112
$ws.Cells.Item(6, 15).Value = 0.281657135515876
$ws.Cells.Item(6, 16).Value = 0.281657135515876
$ws.Cells.Item(6, 17).Value = 647.0555601699235
$ws.Cells.Item(6, 18).Value = 5823.500041529312
$ws.Cells.Item(6, 19).Value = 0.2445215247210608
$ws.Cells.Item(6, 20).Value = 0.2445215247210607

# Row 7: FAPs -> Lamc2 -> Itgb1 -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lamc2"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 6.369908666666666
$ws.Cells.Item(7, 8).Value = 19.109726
$ws.Cells.Item(7, 9).Value = 0.8681531333236113
$ws.Cells.Item(7, 10).Value = 0.8681531333236113
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 142.0267893333333
$ws.Cells.Item(7, 14).Value = 426.080368
$ws.Cells.Item(7, 15).Value = 0.3938062999413425
$ws.Cells.Item(7, 16).Value = 0.3938062999413425
$ws.Cells.Item(7, 17).Value = 904.6976762732409
$ws.Cells.Item(7, 18).Value = 8142.279086459168
$ws.Cells.Item(7, 19).Value = 0.3418841732166544
$ws.Cells.Item(7, 20).Value = 0.3418841732166544

# Row 8: sCs -> Lamc2 -> Itgb1 -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Lamc2"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5252196666666666
$ws.Cells.Item(8, 8).Value = 1.575659
$ws.Cells.Item(8, 9).Value = 0.07158204664470585
$ws.Cells.Item(8, 10).Value = 0.07158204664470584
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 117.044563
$ws.Cells.Item(8, 14).Value = 351.133689
$ws.Cells.Item(8, 15).Value = 0.3245365645427815
$ws.Cells.Item(8, 16).Value = 0.3245365645427815
$ws.Cells.Item(8, 17).Value = 61.47410636400566
$ws.Cells.Item(8, 18).Value = 553.266957276051
$ws.Cells.Item(8, 19).Value = 0.02323099150101397
$ws.Cells.Item(8, 20).Value = 0.02323099150101397

# Row 9: sCs -> Lamc2 -> Itgb1 -> FAPs
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Lamc2"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5252196666666666
$ws.Cells.Item(9, 8).Value = 1.575659
$ws.Cells.Item(9, 9).Value = 0.07158204664470585
$ws.Cells.Item(9, 10).Value = 0.07158204664470584
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 101.5800373333333
$ws.Cells.Item(9, 14).Value = 304.740112
$ws.Cells.Item(9, 15).Value = 0.281657135515876
$ws.Cells.Item(9, 16).Value = 0.281657135515876
$ws.Cells.Item(9, 17).Value = 53.35183334820089
$ws.Cells.Item(9, 18).Value = 480.166500133808
$ws.Cells.Item(9, 19).Value = 0.02016159421231167
$ws.Cells.Item(9, 20).Value = 0.02016159421231167

# Row 10: sCs -> Lamc2 -> Itgb1 -> sCs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Lamc2"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5252196666666666
$ws.Cells.Item(10, 8).Value = 1.575659
$ws.Cells.Item(10, 9).Value = 0.07158204664470585
$ws.Cells.Item(10, 10).Value = 0.07158204664470584
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 142.0267893333333
$ws.Cells.Item(10, 14).Value = 426.080368
$ws.Cells.Item(10, 15).Value = 0.3938062999413425
$ws.Cells.Item(10, 16).Value = 0.3938062999413425
$ws.Cells.Item(10, 17).Value = 74.59526295139023
$ws.Cells.Item(10, 18).Value = 671.357366562512
$ws.Cells.Item(10, 19).Value = 0.02818946093138021
$ws.Cells.Item(10, 20).Value = 0.0281894609313802

